# Weekly refresh: a new day's quote is inserted at row 36 (pushing the
# existing rows 36-136 down to 37-137); the former last row (136) simply
# slides down to 137 with its data unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 36, shifting everything
# below it (through the old row 136) down by one.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 44925
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100102
$ws.Cells.Item(36, 8).Value = "Cítricos"
$ws.Cells.Item(36, 9).Value = 100102004
$ws.Cells.Item(36, 10).Value = "Mandarina"
$ws.Cells.Item(36, 11).Value = "Murcott"
$ws.Cells.Item(36, 12).Value = "Segunda"
$ws.Cells.Item(36, 13).Value = 250
$ws.Cells.Item(36, 14).Value = 15000
$ws.Cells.Item(36, 15).Value = 16000
$ws.Cells.Item(36, 16).Value = 15400
$ws.Cells.Item(36, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(36, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(36, 19).Value = 770
$ws.Cells.Item(36, 20).Value = 20
